# Adapt column header formatting to respective input file names
# - Rename header strings "<name>_old" -> "<name>_FV2304" and "<name>_new" -> "<name>_FV2310"
# - Turn the used range into an Excel Table ("Table1")
# - Freeze the header row (top row)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 21  # column U
$lastRow = 80

for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value2
    if ($val -ne $null) {
        if ($val.EndsWith("_old")) {
            $cell.Value2 = $val.Substring(0, $val.Length - 4) + "_FV2304"
        } elseif ($val.EndsWith("_new")) {
            $cell.Value2 = $val.Substring(0, $val.Length - 4) + "_FV2310"
        }
    }
}

# Create the table over the full used range, using the header row for names
$range = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$tbl = $ws.ListObjects.Add(1, $range, $null, 1)
$tbl.Name = "Table1"

# Freeze the top (header) row
$ws.Application.ActiveWindow.FreezePanes = $false
$ws.Range("A2").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
